# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.921.71"
$ws.Range("E2").Value = "  -1.38%  "

# Row 3
$ws.Range("D3").Value = "3.511.49"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.67"
$ws.Range("E5").Value = "  +4.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.28"
$ws.Range("E6").Value = "  -6.19%  "

# Row 7
$ws.Range("E7").Value = "  +4.66%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  -0.37%  "

# Row 10
$ws.Range("E10").Value = "  +4.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.29"
$ws.Range("E11").Value = "  -0.67%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000272"
$ws.Range("E12").Value = "  +1.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.22"
$ws.Range("E13").Value = "  -2.14%  "

# Row 14
$ws.Range("D14").Value = "4.079.91"
$ws.Range("E14").Value = "  +0.21%  "

# Row 15
$ws.Range("D15").Value = "3.515.97"
$ws.Range("E15").Value = "  +0.27%  "

# Row 16
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.35"
$ws.Range("E17").Value = "  +0.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.07"
$ws.Range("E18").Value = "  +1.89%  "

# Row 19
$ws.Range("D19").Value = "65.927.00"
$ws.Range("E19").Value = "  -1.32%  "

# Row 20
$ws.Range("E20").Value = "  +0.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.68"
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  +8.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  +3.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.83"
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.99"
$ws.Range("E25").Value = "  +8.97%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.01"
$ws.Range("E26").Value = "  -1.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  -2.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("E28").Value = "  +2.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.31"
$ws.Range("E29").Value = "  -0.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "624.66"
$ws.Range("E30").Value = "  -4.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.46"
$ws.Range("E31").Value = "  -4.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.66"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  -0.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  +14.20%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.73"
$ws.Range("E35").Value = "  -0.18%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0800"
$ws.Range("E37").Value = "  -1.34%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.25"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39
$ws.Range("D39").Value = "3.308.66"
$ws.Range("E39").Value = "  +10.21%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("E40").Value = "  +4.02%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.380"
$ws.Range("E41").Value = "  -3.44%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  -0.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0417"
$ws.Range("E44").Value = "  -0.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -5.25%  "

# Row 46
$ws.Range("E46").Value = "  -3.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  +0.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.132"
$ws.Range("E48").Value = "  +1.57%  "

# Row 49
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.57"
$ws.Range("E49").Value = "  -4.47%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.12"
$ws.Range("E50").Value = "  -0.06%  "

# Row 51
$ws.Range("E51").Value = "  -4.28%  "

